$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3..62 down to 4..63
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new record
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 45257
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 100112026
$ws.Range("G3").Value = "Haba"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11500
$ws.Range("N3").Value = '$/saco 25 kilos'
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 460
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
